$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("J1").Value = "incorp"
$ws.Range("K1").Value = "t.incorp"
$ws.Range("L1").Value = "EF.report"

# Apply center alignment to the new K column (t.incorp) for all data rows
$ws.Range("K2:K35").HorizontalAlignment = -4108

# Row 2 (id=1)
$ws.Range("J2").Value = "None"
$ws.Range("L2").Value = 1.6

# Row 3 (id=2)
$ws.Range("J3").Value = "None"
$ws.Range("L3").Value = 1.6

# Row 4 (id=3)
$ws.Range("J4").Value = "None"
$ws.Range("L4").Value = 12

# Row 5 (id=4)
$ws.Range("J5").Value = "None"
$ws.Range("L5").Value = 14

# Row 6 (id=5)
$ws.Range("J6").Value = "None"
$ws.Range("L6").Value = 15

# Row 7 (id=6)
$ws.Range("J7").Value = "None"
$ws.Range("L7").Value = 1.6

# Row 8 (id=7)
$ws.Range("J8").Value = "None"
$ws.Range("L8").Value = 15

# Row 9 (id=8)
$ws.Range("J9").Value = "None"
$ws.Range("L9").Value = 1.6

# Row 10 (id=9)
$ws.Range("J10").Value = "Deep"
$ws.Range("K10").Value = 4
$ws.Range("L10").Value = 3.1

# Row 11 (id=10)
$ws.Range("J11").Value = "Deep"
$ws.Range("K11").Value = 4
$ws.Range("L11").Value = 4.7

# Row 12 (id=11)
$ws.Range("J12").Value = "None"
$ws.Range("L12").Value = 22

# Row 13 (id=12)
$ws.Range("J13").Value = "None"
$ws.Range("L13").Value = 27

# Row 14 (id=13)
$ws.Range("J14").Value = "None"
$ws.Range("L14").Value = 30

# Row 15 (id=14)
$ws.Range("F15").Value = 16.86667
$ws.Range("G15").Value = 3.1816667
$ws.Range("J15").Value = "None"
$ws.Range("L15").Value = 33

# Row 16 (id=15)
$ws.Range("F16").Value = 16.86667
$ws.Range("G16").Value = 3.1816667
$ws.Range("J16").Value = "Deep"
$ws.Range("K16").Value = 4
$ws.Range("L16").Value = 11

# Row 17 (id=16)
$ws.Range("J17").Value = "None"
$ws.Range("L17").Value = 32

# Row 18 (id=17)
$ws.Range("J18").Value = "Deep"
$ws.Range("K18").Value = 4
$ws.Range("L18").Value = 9.2

# Row 19 (id=18)
$ws.Range("J19").Value = "None"
$ws.Range("L19").Value = 1.8

# Row 20 (id=19)
$ws.Range("J20").Value = "None"
$ws.Range("L20").Value = 1.8

# Row 21 (id=20)
$ws.Range("J21").Value = "None"
$ws.Range("L21").Value = 9.7

# Row 22 (id=21)
$ws.Range("J22").Value = "None"
$ws.Range("L22").Value = 9.8

# Row 23 (id=22)
$ws.Range("J23").Value = "None"
$ws.Range("L23").Value = 10

# Row 24 (id=23)
$ws.Range("J24").Value = "None"
$ws.Range("L24").Value = 1.8

# Row 25 (id=24)
$ws.Range("J25").Value = "None"
$ws.Range("L25").Value = 10

# Row 26 (id=25)
$ws.Range("J26").Value = "None"
$ws.Range("L26").Value = 1.8

# Row 27 (id=26)
$ws.Range("J27").Value = "Deep"
$ws.Range("K27").Value = 4
$ws.Range("L27").Value = 2.4

# Row 28 (id=27)
$ws.Range("J28").Value = "Deep"
$ws.Range("K28").Value = 4
$ws.Range("L28").Value = 3.5

# Row 29 (id=28)
$ws.Range("J29").Value = "None"
$ws.Range("L29").Value = 14

# Row 30 (id=29)
$ws.Range("J30").Value = "None"
$ws.Range("L30").Value = 15

# Row 31 (id=30)
$ws.Range("J31").Value = "None"
$ws.Range("L31").Value = 15

# Row 32 (id=31)
$ws.Range("F32").Value = 16.86667
$ws.Range("G32").Value = 3.1816667
$ws.Range("J32").Value = "None"
$ws.Range("L32").Value = 16

# Row 33 (id=32)
$ws.Range("F33").Value = 16.86667
$ws.Range("G33").Value = 3.1816667
$ws.Range("J33").Value = "Deep"
$ws.Range("K33").Value = 4
$ws.Range("L33").Value = 6.7

# Row 34 (id=33)
$ws.Range("J34").Value = "None"
$ws.Range("L34").Value = 15

# Row 35 (id=34)
$ws.Range("J35").Value = "Deep"
$ws.Range("K35").Value = 4
$ws.Range("L35").Value = 5.8
